# Applies the Review_282.docx edit: refreshed date, new paper title,
# rewritten body paragraphs, and the body paragraph that used to hold
# the arXiv link is expanded into four paragraphs ending with the new link.
$d = $word.ActiveDocument

$replacements = @(
    @("⚡️🚀המאמר היומי של מייק 22.08.24: ⚡️🚀", "⚡️🚀המאמר היומי של מייק 21.08.24: ⚡️🚀"),
    @("Approaching Deep Learning through the Spectral Dynamics of Weights", "Tree Attention: Topology-Aware Decoding for Long-Context Attention on GPU Clusters"),
    @("היום נסקור מאמר החוקר מה הסיבות לתופעה של גרוקינג. למי שלא מכיר גרוקינג זו תופעה די מעניינת המתרחשת כאשר ממשיכים לאמן רשת נוירונים (למרות שזה קורה גם במודלים אחרים) גם אחרי לוס הוולידציה מתחיל לעלות (כלומר אנו נכנסים למשטר אוורפיט). מתברר אם לא עוצרים וממשיכים לאמן לוס הוולידציה מתחיל לרדת כלומר המודל נכנס למשטר ההכללה כלומר לומד את ה״חוקיות האמיתית״ מאחורי הדאטה. ", "היום נסקור מאמר בנושאה שכבר סקרתי כמה מאמרים לפני כחודש. הנושא הזה נקרא אופטימיזציה והאצה decoding של מודלי שפה כלומר התהליך שגנרוט טוקן חדש בתלות בכל הטוקנים בתוך חלון ההקשר שכבר גונרטו. ואם חלון ההקשר הוא ארוך (מאות אלפי טוקנים) זה יכול לקחת די הרבה זמן בעיקר בגלל מנגנון ה-attention של הטרנספורמרים שמהווים backbone של כל מודלי השפה החזקים."),
    @("התופעה הזו היא מקרה פרטי של double descent (יש גם multiple descent) שמתרחש גם אם אנו מוסיפים פרמטרים למודל בצורה עקבית ומגיעים למצב שיש לנו over-parametrization. כלומר יש המודל שלנו לכאורה מתחיל ״יותר מדי פרמטרים״ כדי ״להבין את הדאטה״. וגם שם זה קורה בצורה בלתי רציפה כלומר יש אינטרוול של פרמטרים שביצועי המודל יורדים עבורם ורק אז מתחילים לרדת. ", "בשנים האחרונות הוצעו מספר רב של שיטות לייעול והאצה של חישוב ה-attention שהכי מפורסמים מהם הם Flash Attention ו-KV-Cache. שיטות אלו בדרך כלל מנצלות את העובדה שהיום אינפרנס של מודלי שפה מתבצע על GPU וניתן לייעל את החישוב על ידי שימוש ביכולת של GPUs לחשב דברים במקביל."),
    @("המאמר חוקר מה קורה עם משקלי המודל כאשר הוא נכנס למשטר הגרוקינג. מתברר שתופעה הגרוקינג קשורה לירידה בראנק של מטריצות המשקלים של המודל. בשבילי זה די אינטואיטיבי כי לדעתי במהלך גרוקינג המודל מצליח להתכנס ל״פתרון פשוט ביותר עבור הדאטהסט. פתרון פשוט הכוונה הוא מודל שאפקטיבית הוא קטן, כלומר רוב וקטורי המשקלים בו או אפס או תלוים לינארית זה בזה.", "יתרה מזו מכיוון שמודלי שפה רצים היום על קלסטרים של GPUs יצאו מספר עבודות על איך ניתן לחשב את ה-attention על קלסטרים אלו. מכיוון שמנגנון ה-attention מכיל מכפלות פנימיות (סכומים רבים) אז ניתן לחשבו בצורה מבוזרת די ביעילות. ")
)

foreach ($pair in $replacements) {
    $ok = $d.Content.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)
    if (-not $ok) { throw "replacement not found: $($pair[0])" }
}

# The last paragraph (the arXiv link) grows into four paragraphs: three new
# discussion paragraphs followed by the (changed) link. InsertBefore pushes
# each new paragraph immediately above the link paragraph, so insert them
# in reverse order to land in the right reading order.
$pLink = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLink.Range.InsertBefore("למה זה טוב בכלל? מתברר שהחישוב של attention בצורה כזו מערב פעולות כמו logsumexp ו- max שניתן לבזר אותם בצורה יעילה בין ה-GPUs. החישוב נעשה בצורה של עץ, כלומר מחלקים את הסכומים לכמה חלקים, מחשבים כל חלק ואז מתחילים לסכם את התוצאות בצורה היררכית. זה כמו Map-Reduce רב שלבי.  `r")
$pLink.Range.InsertBefore("ניתן להכליל את החישוב הזה ל-attention עבור וקטורי שאילתה q מרובים כאשר במקום נגזרת רגילה יהיה לנו נגזרת לפי n משתנים (n הינו מספר וקטורי השאילתה). `r")
$pLink.Range.InsertBefore("והמאמר הזה מציע מנגנון מעניין של חישוב ה-attention. הדבר המעניין בו שהמאמר הזה מייצג את חישוב ה-attention (עבור וקטור שאילתה נתון q) כנגזרת של הלוג של ״פונקציה יוצרת״ של ה-attention המחושבת בנקודת 0. פונקציה יוצרת זו נבנית על ידי מניפולציה פשוטה של נוסחת ה-attention וממש מזכירה פונקציה יוצרת של משתנה אקראי. `r")

$okLink = $d.Content.Find.Execute("https://arxiv.org/abs/2408.11804", $false, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2408.04093", 2)
if (-not $okLink) { throw "link replacement not found" }

Write-Output "paragraphs: $($d.Paragraphs.Count)"
